$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing row (69) onto the new row 70,
# so the new date cell picks up the same date-formatted style (s="1").
$ws.Range("A69").Copy()
$ws.Range("A70").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row's values (row 70: 2024-09-01 data point)
$ws.Range("A70").Value = 45536
$ws.Range("B70").Value = -0.274
$ws.Range("C70").Value = -0.646
$ws.Range("D70").Value = 0.177
$ws.Range("E70").Value = -0.109
$ws.Range("F70").Value = -0.906
$ws.Range("G70").Value = 2.1
